$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" text cell (A1) with the new conversion rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.86 = 10673.77 pesos`n✅ 10673.77 pesos = 2.84 = 951.9 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the "tasas" sheet numeric values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 349.777
$ws2.Range("O10").Value = 3733.44
$ws2.Range("N12").Value = 3760
$ws2.Range("O12").Value = 335.322
